$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-15 Monday" "2024-07-16 Tuesday"

Replace-Text "81÷5=" "81÷7="
Replace-Text "31÷6=" "59÷5="
Replace-Text "81÷6=" "43÷5="
Replace-Text "17÷6=" "34÷2="
Replace-Text "12÷7=" "77÷9="

# Note: "58÷6=" must be replaced before "80÷6=" is turned into "58÷6=",
# otherwise the replace-all for 80÷6= would create a duplicate that a
# later pass could mistakenly touch.
Replace-Text "58÷6=" "18÷3="
Replace-Text "80÷6=" "58÷6="

Replace-Text "18÷6=" "63÷9="
Replace-Text "39÷3=" "99÷6="
Replace-Text "36÷3=" "95÷4="
Replace-Text "11÷5=" "72÷7="
Replace-Text "20÷2=" "62÷6="
Replace-Text "26÷2=" "98÷5="
Replace-Text "16÷7=" "27÷2="
Replace-Text "80÷3=" "43÷9="
Replace-Text "26÷8=" "51÷5="
Replace-Text "71÷7=" "92÷4="
Replace-Text "61÷3=" "54÷4="
Replace-Text "69÷9=" "11÷3="
Replace-Text "65÷8=" "82÷2="
Replace-Text "65÷3=" "59÷3="
Replace-Text "25÷5=" "22÷3="
Replace-Text "64÷2=" "55÷8="
Replace-Text "56÷7=" "85÷8="
Replace-Text "43÷4=" "71÷6="
